$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1) updates
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 7070
$ws1.Range("F5").Value = 461
$ws1.Range("F7").Value = 7627
$ws1.Range("F13").Value = 436
$ws1.Range("F14").Value = 161
$ws1.Range("F17").Value = 57
$ws1.Range("F18").Value = 57
$ws1.Range("F19").Value = 24
$ws1.Range("F20").Value = 5482
$ws1.Range("F21").Value = 142
$ws1.Range("F22").Value = 201
$ws1.Range("F23").Value = 913
$ws1.Range("F25").Value = 298

# Sheet "全部类型" (sheet4) updates
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 7070
$ws4.Range("F5").Value = 461
$ws4.Range("F7").Value = 7627
$ws4.Range("F13").Value = 436
$ws4.Range("F14").Value = 161
$ws4.Range("F17").Value = 57
$ws4.Range("F18").Value = 57
$ws4.Range("F19").Value = 24
$ws4.Range("F21").Value = 5482
$ws4.Range("F23").Value = 142
$ws4.Range("F24").Value = 201
$ws4.Range("F25").Value = 913
$ws4.Range("F27").Value = 298
